# Add "2022-Q1" sheet (copied structure/style from "2021-Q4") with new fund data,
# positioned right before "总计", and add a corresponding new row to "总计".

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q1" sheet by copying "2021-Q4" (same columns/styles) ---
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$q4Sheet.Copy($null, $q4Sheet)
$newSheet = $wb.Worksheets.Item($q4Sheet.Index + 1)
$newSheet.Name = "2022-Q1"

# Helper: a cell with no explicit style to use as a formatting donor for text cells,
# so the text keeps numFmt/quotePrefix-free styling identical to the header copy.
$plainFormatCell = $newSheet.Cells.Item(2, 2)

function Set-TextValue($cell, [string]$text) {
    $cell.Value = "'" + $text
    $plainFormatCell.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}

# Row 2: fund 004845 (南华瑞盈混合A) -- code/name unchanged, other metrics updated
Set-TextValue $newSheet.Cells.Item(2, 4) "2.90"
Set-TextValue $newSheet.Cells.Item(2, 5) "82.72"
Set-TextValue $newSheet.Cells.Item(2, 6) "2.81"
Set-TextValue $newSheet.Cells.Item(2, 7) "0.0815"
$newSheet.Cells.Item(2, 8).Value = 5

# Row 3: fund 004846 (南华瑞盈混合C) -- code/name unchanged, other metrics updated
Set-TextValue $newSheet.Cells.Item(3, 4) "0.10"
Set-TextValue $newSheet.Cells.Item(3, 5) "82.72"
Set-TextValue $newSheet.Cells.Item(3, 6) "2.81"
Set-TextValue $newSheet.Cells.Item(3, 7) "0.0028"
$newSheet.Cells.Item(3, 8).Value = 5

$excel.CutCopyMode = 0

# --- 2. Insert a new top data row into "总计" for the 2022-Q1 summary ---
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Rows.Item(2).ClearFormats()

# Restore the row-index cell's style (A column keeps style used by the other rows)
$totalSheet.Cells.Item(3, 1).Copy() | Out-Null
$totalSheet.Cells.Item(2, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 2
$totalSheet.Cells.Item(2, 4).Value = 0.08

# Re-number the row-index column (A) for the rows that shifted down
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
